# Apply the "through August 12" data refresh to the carjacking-by-neighborhood
# workbook:
#   - bump the as-of date in the sheet name and the first data-column header
#   - add a handful of newly-geocoded incidents to existing neighborhoods
#     (all land in an "August <year>" column)
#   - insert a brand-new neighborhood row ("Magnificent Mile") in its correct
#     alphabetical slot, with its own new incident
#   - one more new incident for the last row (West Ridge) in its August 2020 column

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sheet tab + update the "as of" header text -------------------
$ws.Name = "Through 2022-08-12"
$ws.Range("B1").Value = "August 2022 (through August 12)"

# --- Scattered incremental counts on existing neighborhoods ---------------
$ws.Range("B2").Value = 10
$ws.Range("J2").Value = 6
$ws.Range("Z3").Value = 3
$ws.Range("AX3").Value = 2
$ws.Range("B5").Value = 2
$ws.Range("AP8").Value = 3
$ws.Range("AX9").Value = 4
$ws.Range("J11").Value = 2
$ws.Range("AH13").Value = 3
$ws.Range("AP17").Value = 1
$ws.Range("J28").Value = 2
$ws.Range("R28").Value = 1
$ws.Range("BF28").Value = 1
$ws.Range("AX29").Value = 1
$ws.Range("B31").Value = 1
$ws.Range("R41").Value = 1
$ws.Range("B42").Value = 1

# --- Insert the new "Magnificent Mile" row (alphabetically after --------
# --- "Lincoln Square", before "Mckinley Park", i.e. row 79) --------------
$ws.Rows(79).Insert()

# match the formatting (bold/centered/top/bordered label style) used by
# every other neighborhood-name cell in column A
$ws.Range("A80").Copy()
$ws.Range("A79").PasteSpecial(-4122)

$ws.Range("A79").Value = "Magnificent Mile"
$ws.Range("AH79").Value = 1

# --- One more new incident on the last row (West Ridge, now row 98) ------
$ws.Range("R98").Value = 2
